# Update "想去人数" (interest count) and "最低票价" (min ticket price) figures
# for the refreshed data pull (gh-pages output regenerated at 456a3b4).
#
# The workbook has 4 sheets:
#   1 - 展览     (Exhibitions)      -> needs updates
#   2 - 演出     (Performances)     -> unchanged
#   3 - 本地生活 (Local life)       -> unchanged (empty)
#   4 - 全部类型 (All types)        -> needs updates (same rows as sheet 1,
#                                     shifted by one extra row because it
#                                     also contains the performance entry)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 --------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value  = 6658
$ws1.Range("F3").Value  = 9
$ws1.Range("F6").Value  = 10
$ws1.Range("F7").Value  = 544
$ws1.Range("F11").Value = 4
$ws1.Range("F14").Value = 1289
$ws1.Range("F16").Value = 3320
$ws1.Range("F18").Value = 214
$ws1.Range("F19").Value = 1968
$ws1.Range("F20").Value = 80
$ws1.Range("G20").Value = 39.9

# --- Sheet 4: 全部类型 -----------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F2").Value  = 6658
$ws4.Range("F3").Value  = 9
$ws4.Range("F6").Value  = 10
$ws4.Range("F8").Value  = 544
$ws4.Range("F12").Value = 4
$ws4.Range("F15").Value = 1289
$ws4.Range("F17").Value = 3320
$ws4.Range("F19").Value = 214
$ws4.Range("F20").Value = 1968
$ws4.Range("F21").Value = 80
$ws4.Range("G21").Value = 39.9
